# Add the "ODI Batting Extra" sheet (sheet3) to the workbook, after the
# existing "ODI Batting" sheet, and populate it with the player's extra
# ODI batting stats (4s, 6s, % of team runs, man-of-the-match flag).

$wb = $excel.ActiveWorkbook

$playerInfo = $wb.Worksheets.Item(1)
$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Match the page margins used by the workbook's other sheets (0.75in
# sides, 1in top/bottom, 0.5in header/footer == 54/54/72/72/36/36 pt).
$ws.PageSetup.LeftMargin   = 54
$ws.PageSetup.RightMargin  = 54
$ws.PageSetup.TopMargin    = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ---- header row -----------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")

# Copy the bold/centered/bordered header style already used by the other
# sheets onto row 1 of the new sheet before filling in the text.
$playerInfo.Range("A1:D1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
}

# ---- data rows --------------------------------------------------------
# Columns: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("3996", 2,    "0", "0", "",       "NO"),
    @("4092", 2,    "0", "1", "4.20%",  "NO"),
    @("4093", 2,    "1", "1", "9.09%",  "NO"),
    @("4129", 2,    "9", "0", "16.22%", "NO"),
    @("4131", 2,    "0", "0", "1.12%",  "NO"),
    @("4132", $null, $null, $null, $null, "NO"),
    @("4134", 2,    "6", "2", "37.78%", "NO"),
    @("4140", 2,    "0", "0", "",       "NO"),
    @("4145", 2,    "0", "0", "0.52%",  "NO"),
    @("4147", $null, $null, $null, $null, "NO"),
    @("4189", 2,    "1", "0", "3.96%",  "NO"),
    @("4192", 2,    "6", "2", "44.88%", "NO"),
    @("4195", 2,    "6", "0", "18.07%", "NO"),
    @("4198", 2,    "2", "0", "3.14%",  "NO"),
    @("4200", 2,    "1", "0", "3.89%",  "NO"),
    @("4202", 2,    "2", "0", "3.25%",  "NO")
)

# Write a value as plain *text* (matching the source data's inline-string
# cells), even when it looks like a number (e.g. "3996", "0"), without
# leaving any lingering number-format style behind on the cell. A quoted
# text formula always evaluates to a string; copy/pasting it back as
# values bakes in that literal text with the sheet's default (unstyled)
# cell format - same as the plain text cells Excel would never have
# coerced to a number in the first place.
function Set-TextCell($cell, $value) {
    if ($value -eq $null -or $value -eq "") {
        $cell.Value = ""
        return
    }
    $escaped = [string]$value -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$r = 2
foreach ($row in $rows) {

    # A: MATCH_CODE (text)
    Set-TextCell $ws.Cells.Item($r, 1) $row[0]

    # B: BATTING_POSITION (number, or blank)
    $cellB = $ws.Cells.Item($r, 2)
    if ($row[1] -eq $null) {
        $cellB.Value = ""
    } else {
        $cellB.Value = $row[1]
    }

    # C: NUM_4 (text, or blank)
    Set-TextCell $ws.Cells.Item($r, 3) $row[2]

    # D: NUM_6 (text, or blank)
    Set-TextCell $ws.Cells.Item($r, 4) $row[3]

    # E: PERCENT_RUNS_OF_TOTAL (text, or blank)
    Set-TextCell $ws.Cells.Item($r, 5) $row[4]

    # F: MAN_OF_MATCH (text)
    Set-TextCell $ws.Cells.Item($r, 6) $row[5]

    $r++
}

# Keep the first sheet active/selected, as it was before the edit.
$playerInfo.Activate()
$playerInfo.Range("A1").Select()
